$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.655.94"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "2.234.58"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  +0.42%  "
$ws.Range("D5").Value = "'305.63"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "'94.48"
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("D7").Value = "'0.571"
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("D9").Value = "'0.516"
$ws.Range("E9").Value = "  -2.08%  "
$ws.Range("D10").Value = "'34.80"
$ws.Range("E10").Value = "  -0.95%  "
$ws.Range("D11").Value = "'0.0803"
$ws.Range("E11").Value = "  -1.52%  "
$ws.Range("D12").Value = "'7.16"
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("D13").Value = "'0.104"
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("D14").Value = "2.575.52"
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("D15").Value = "2.226.55"
$ws.Range("E15").Value = "  -4.54%  "
$ws.Range("D16").Value = "'0.832"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "'13.52"
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").Value = "44.578.16"
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("D19").Value = "0.0₃0943"
$ws.Range("E19").Value = "  -2.75%  "
$ws.Range("D20").Value = "'11.89"
$ws.Range("E20").Value = "  -1.88%  "
$ws.Range("D21").Value = "'6.21"
$ws.Range("E21").Value = "  -3.10%  "
$ws.Range("D22").Value = "'64.98"
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("D23").Value = "'238.47"
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("D24").Value = "'2.92"
$ws.Range("E24").Value = "  -1.06%  "
$ws.Range("E25").Value = "  -1.41%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'2.30"
$ws.Range("E27").Value = "  +3.60%  "
$ws.Range("D28").Value = "'9.72"
$ws.Range("E28").Value = "  -2.59%  "
$ws.Range("D29").Value = "'37.21"
$ws.Range("E29").Value = "  -0.83%  "
$ws.Range("D30").Value = "'5.89"
$ws.Range("E30").Value = "  -1.84%  "
$ws.Range("D31").Value = "'19.83"
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("D32").Value = "'150.17"
$ws.Range("E32").Value = "  -2.08%  "
$ws.Range("D33").Value = "'0.0786"
$ws.Range("E33").Value = "  -1.87%  "
$ws.Range("E34").Value = "  +0.86%  "
$ws.Range("D35").Value = "'3.04"
$ws.Range("E35").Value = "  -6.64%  "
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("E37").Value = "  -2.10%  "
$ws.Range("E38").Value = "  +5.31%  "
$ws.Range("E39").Value = "  +3.77%  "
$ws.Range("D40").Value = "'3.34"
$ws.Range("E40").Value = "  -3.54%  "
$ws.Range("D41").Value = "'3.75"
$ws.Range("E41").Value = "  -3.00%  "
$ws.Range("D42").Value = "'0.0300"
$ws.Range("E42").Value = "  +0.49%  "
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").Value = "1.826.14"
$ws.Range("E44").Value = "  +4.85%  "
$ws.Range("D45").Value = "'1.72"
$ws.Range("E45").Value = "  +10.91%  "
$ws.Range("D46").Value = "'79.51"
$ws.Range("E46").Value = "  -4.36%  "
$ws.Range("D47").Value = "'0.187"
$ws.Range("E47").Value = "  -1.97%  "
$ws.Range("D48").Value = "'98.19"
$ws.Range("E48").Value = "  -1.95%  "
$ws.Range("D49").Value = "'4.85"
$ws.Range("E49").Value = "  -2.14%  "
$ws.Range("D50").Value = "'68.95"
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("D51").Value = "'7.97"
$ws.Range("E51").Value = "  -2.40%  "
